$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1 / rId1) - F column updates
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 1585
$wsExpo.Range("F5").Value = 11
$wsExpo.Range("F8").Value = 138
$wsExpo.Range("F9").Value = 62
$wsExpo.Range("F10").Value = 458

# Sheet "全部类型" (index 4 / rId4) - F column updates
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 107
$wsAll.Range("F4").Value = 1585
$wsAll.Range("F6").Value = 23
$wsAll.Range("F9").Value = 62
$wsAll.Range("F10").Value = 458
